$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 3075
$wsExpo.Range("F4").Value = 144
$wsExpo.Range("F5").Value = 98

# Sheet "演出" (Performances)
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 119

# Sheet "全部类型" (All types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 119
$wsAll.Range("F7").Value = 3075
$wsAll.Range("F8").Value = 144
$wsAll.Range("F10").Value = 98
